# Insert a new slide ("Model Execution Mode") right after "Technology Stack"
# (slide 12) and before "Experimental Strategy" (old slide 13). Every slide
# from the old position 13 onward is pushed down by one, ending with a
# 23-slide deck whose final slide is still "Thank You".

$p = $ppt.ActivePresentation

# Slide 13 currently holds "Experimental Strategy" - insert the new slide
# before it, reusing the same "Title and Content" layout (as used by slide 13).
$layout = $p.Slides.Item(13).Layout
$s = $p.Slides.Add(13, $layout)

# Title placeholder
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Model Execution Mode"

# Body / content placeholder - four bullet paragraphs, 22pt each (matches the
# sz="2200" styling used throughout the rest of the deck's bullet slides).
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Current implementation uses local models only`r" + `
             "No OpenAI API key is required for this submission`r" + `
             "Offline mock mode is available for restricted environments`r" + `
             "OpenAI integration can be added later as future scope"
$body.Font.Size = 22

Write-Output "Inserted slide 13 'Model Execution Mode'; deck now has $($p.Slides.Count) slides."
